# Actualizacion final con metrica F1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New feature importances (retrained model, F1 metric)
$features = @(
    @{ Name = "Age";                      Importance = 0.1547037571708041 },
    @{ Name = "DistanceFromHome";         Importance = 0.1256262200658769 },
    @{ Name = "EnvironmentSatisfaction";  Importance = 0.07180021619067051 },
    @{ Name = "MonthlyIncome";            Importance = 0.231762606920802 },
    @{ Name = "NumCompaniesWorked";       Importance = 0.08045655969826347 },
    @{ Name = "PercentSalaryHike";        Importance = 0.1078954508922652 },
    @{ Name = "TotalWorkingYears";        Importance = 0.1317218604522631 },
    @{ Name = "YearsAtCompany";           Importance = 0.09603332860905472 }
)

$rowCount = $features.Count

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $features[$i].Name
    $ws.Cells.Item($row, 3).Value = $features[$i].Importance
}

# The new rows (8 and 9) were added beyond the original A1:C7 range; make
# sure column A picks up the same formatting (border/bold/center) used by
# the rest of the index column.
$ws.Range("A7").Copy() | Out-Null
$ws.Range("A8:A9").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
